# Shift the "date" column (A2:A49) forward by one week (7 days),
# keeping the values stored as plain text strings (e.g. "2023-09-18" -> "2023-09-25"),
# and move the active selection from B47 to B45.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 49
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $val = $cell.Value2
    if ($val) {
        $d = [DateTime]::ParseExact($val, "yyyy-MM-dd", [System.Globalization.CultureInfo]::InvariantCulture)
        $newD = $d.AddDays(7)
        $cell.Value2 = $newD.ToString("yyyy-MM-dd")
    }
}

$ws.Range("B45").Select()
